$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: add the two new trailing columns (PriceChange / UpDown) ---
$ws.Cells.Item(5,24).Value = 0.18999999999999773
$ws.Cells.Item(5,25).Value = "Up"

# --- Copy row 5's cell formatting down into the new rows 6-11 ---
# Columns A:W for every new row...
$ws.Range("A5:W5").Copy()
$ws.Range("A6:W11").PasteSpecial(-4122)
# ...and columns X:Y only for rows 6-10 (row 11 has no PriceChange/UpDown values)
$ws.Range("X5:Y5").Copy()
$ws.Range("X6:Y10").PasteSpecial(-4122)

# --- Row 6 ---
$ws.Cells.Item(6,1).Value = 42650.338368055556
$ws.Cells.Item(6,2).Value = 11
$ws.Cells.Item(6,3).Value = "Buy"
$ws.Cells.Item(6,4).Value = 50
$ws.Cells.Item(6,5).Value = 6515
$ws.Cells.Item(6,6).Value = 354
$ws.Cells.Item(6,7).Value = 67
$ws.Cells.Item(6,8).Value = 32
$ws.Cells.Item(6,9).Value = 91
$ws.Cells.Item(6,10).Value = 8
$ws.Cells.Item(6,11).Value = 10648
$ws.Cells.Item(6,12).Value = 83
$ws.Cells.Item(6,13).Value = 40
$ws.Cells.Item(6,14).Value = 34
$ws.Cells.Item(6,15).Value = 3
$ws.Cells.Item(6,16).Value = "Named"
$ws.Cells.Item(6,17).Value = 38.48959524716075
$ws.Cells.Item(6,18).Value = 0
$ws.Cells.Item(6,19).Value = 0.1046
$ws.Cells.Item(6,20).Value = 0.0345
$ws.Cells.Item(6,21).Value = 4.82
$ws.Cells.Item(6,22).Value = 2.28
$ws.Cells.Item(6,23).Value = 0
$ws.Cells.Item(6,24).Value = 0.18999999999999773
$ws.Cells.Item(6,25).Value = "Up"

# --- Row 7 ---
$ws.Cells.Item(7,1).Value = 42650.339606481481
$ws.Cells.Item(7,2).Value = -6
$ws.Cells.Item(7,3).Value = "Neutral"
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(7,5).Value = 0
$ws.Cells.Item(7,6).Value = 0
$ws.Cells.Item(7,7).Value = 0
$ws.Cells.Item(7,8).Value = 0
$ws.Cells.Item(7,9).Value = 0
$ws.Cells.Item(7,10).Value = 0
$ws.Cells.Item(7,11).Value = 2025
$ws.Cells.Item(7,12).Value = 0
$ws.Cells.Item(7,13).Value = 0
$ws.Cells.Item(7,14).Value = 0
$ws.Cells.Item(7,15).Value = 0
$ws.Cells.Item(7,16).Value = "Named"
$ws.Cells.Item(7,17).Value = 38.48959524716075
$ws.Cells.Item(7,18).Value = 0
$ws.Cells.Item(7,19).Value = 0.1046
$ws.Cells.Item(7,20).Value = 0.0345
$ws.Cells.Item(7,21).Value = 4.82
$ws.Cells.Item(7,22).Value = 2.28
$ws.Cells.Item(7,23).Value = 0
$ws.Cells.Item(7,24).Value = 0.18999999999999773
$ws.Cells.Item(7,25).Value = "Up"

# --- Row 8 ---
$ws.Cells.Item(8,1).Value = 42650.348773148151
$ws.Cells.Item(8,2).Value = 11
$ws.Cells.Item(8,3).Value = "Buy"
$ws.Cells.Item(8,4).Value = 34
$ws.Cells.Item(8,5).Value = 26524
$ws.Cells.Item(8,6).Value = 1445
$ws.Cells.Item(8,7).Value = 61
$ws.Cells.Item(8,8).Value = 35
$ws.Cells.Item(8,9).Value = 88
$ws.Cells.Item(8,10).Value = 12
$ws.Cells.Item(8,11).Value = 32349
$ws.Cells.Item(8,12).Value = 294
$ws.Cells.Item(8,13).Value = 170
$ws.Cells.Item(8,14).Value = 44
$ws.Cells.Item(8,15).Value = 6
$ws.Cells.Item(8,16).Value = "Named"
$ws.Cells.Item(8,17).Value = 38.48959524716075
$ws.Cells.Item(8,18).Value = 0
$ws.Cells.Item(8,19).Value = 0.1046
$ws.Cells.Item(8,20).Value = 0.0345
$ws.Cells.Item(8,21).Value = 4.82
$ws.Cells.Item(8,22).Value = 2.28
$ws.Cells.Item(8,23).Value = 0
$ws.Cells.Item(8,24).Value = 0.18999999999999773
$ws.Cells.Item(8,25).Value = "Up"

# --- Row 9 ---
$ws.Cells.Item(9,1).Value = 42650.359039351853
$ws.Cells.Item(9,2).Value = -4
$ws.Cells.Item(9,3).Value = "Neutral"
$ws.Cells.Item(9,4).Value = -24
$ws.Cells.Item(9,5).Value = 2380
$ws.Cells.Item(9,6).Value = 138
$ws.Cells.Item(9,7).Value = 51
$ws.Cells.Item(9,8).Value = 48
$ws.Cells.Item(9,9).Value = 0
$ws.Cells.Item(9,10).Value = 100
$ws.Cells.Item(9,11).Value = 5395
$ws.Cells.Item(9,12).Value = 18
$ws.Cells.Item(9,13).Value = 17
$ws.Cells.Item(9,14).Value = 0
$ws.Cells.Item(9,15).Value = 1
$ws.Cells.Item(9,16).Value = "Named"
$ws.Cells.Item(9,17).Value = 38.48959524716075
$ws.Cells.Item(9,18).Value = 0
$ws.Cells.Item(9,19).Value = 0.1046
$ws.Cells.Item(9,20).Value = 0.0345
$ws.Cells.Item(9,21).Value = 4.82
$ws.Cells.Item(9,22).Value = 2.28
$ws.Cells.Item(9,23).Value = 0
$ws.Cells.Item(9,24).Value = 0.18999999999999773
$ws.Cells.Item(9,25).Value = "Up"

# --- Row 10 ---
$ws.Cells.Item(10,1).Value = 42650.36146990741
$ws.Cells.Item(10,2).Value = 11
$ws.Cells.Item(10,3).Value = "Buy"
$ws.Cells.Item(10,4).Value = 50
$ws.Cells.Item(10,5).Value = 6500
$ws.Cells.Item(10,6).Value = 352
$ws.Cells.Item(10,7).Value = 67
$ws.Cells.Item(10,8).Value = 32
$ws.Cells.Item(10,9).Value = 91
$ws.Cells.Item(10,10).Value = 8
$ws.Cells.Item(10,11).Value = 10289
$ws.Cells.Item(10,12).Value = 83
$ws.Cells.Item(10,13).Value = 40
$ws.Cells.Item(10,14).Value = 34
$ws.Cells.Item(10,15).Value = 3
$ws.Cells.Item(10,16).Value = "Named"
$ws.Cells.Item(10,17).Value = 38.48959524716075
$ws.Cells.Item(10,18).Value = 0
$ws.Cells.Item(10,19).Value = 0.1046
$ws.Cells.Item(10,20).Value = 0.0345
$ws.Cells.Item(10,21).Value = 4.82
$ws.Cells.Item(10,22).Value = 2.28
$ws.Cells.Item(10,23).Value = 0
$ws.Cells.Item(10,24).Value = 0.18999999999999773
$ws.Cells.Item(10,25).Value = "Up"

# --- Row 11 (no PriceChange / UpDown columns) ---
$ws.Cells.Item(11,1).Value = 42650.36309027778
$ws.Cells.Item(11,2).Value = 1
$ws.Cells.Item(11,3).Value = "Neutral"
$ws.Cells.Item(11,4).Value = 0
$ws.Cells.Item(11,5).Value = 0
$ws.Cells.Item(11,6).Value = 0
$ws.Cells.Item(11,7).Value = 0
$ws.Cells.Item(11,8).Value = 0
$ws.Cells.Item(11,9).Value = 0
$ws.Cells.Item(11,10).Value = 0
$ws.Cells.Item(11,11).Value = 2037
$ws.Cells.Item(11,12).Value = 0
$ws.Cells.Item(11,13).Value = 0
$ws.Cells.Item(11,14).Value = 0
$ws.Cells.Item(11,15).Value = 0
$ws.Cells.Item(11,16).Value = "Named"
$ws.Cells.Item(11,17).Value = 37.799019424898844
$ws.Cells.Item(11,18).Value = 0
$ws.Cells.Item(11,19).Value = 0.1046
$ws.Cells.Item(11,20).Value = 0.0343
$ws.Cells.Item(11,21).Value = 4.82
$ws.Cells.Item(11,22).Value = 2.28
$ws.Cells.Item(11,23).Value = 0

# --- Active cell / selection moves to B7, matching the saved view state ---
$ws.Range("B7").Select()
